$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# 1) Row "01" (table row 2): date "25/10" -> "01/11"
# ---------------------------------------------------------------------------
$cellDate1 = $t.Cell(2, 2)
$cellDate1.Range.Text = "01/11"

# ---------------------------------------------------------------------------
# 2) Row "01" description: "...que serão realizadas na Sprint-7" ->
#    "...que serão realizadas na Sprint-8"
# ---------------------------------------------------------------------------
$cellDesc1 = $t.Cell(2, 4)
$desc1Text = $cellDesc1.Range.Text
$desc1Start = $cellDesc1.Range.Start
$oldStr1 = "que serão realizadas na Sprint-7"
$newStr1 = "que serão realizadas na Sprint-8"
$idx1 = $desc1Text.IndexOf($oldStr1)
if ($idx1 -ge 0) {
    $rngStart = $desc1Start + $idx1
    $rngEnd = $rngStart + $oldStr1.Length
    $d.Range($rngStart, $rngEnd).Text = $newStr1
}

# ---------------------------------------------------------------------------
# 3) Row "03" description: "Retrospectiva da Sprint-7 e início da Sprint-8"
#    -> "Retrospectiva da Sprint-8 e início da Sprint-9"
#    (done before the row above it is deleted, while indices are stable)
# ---------------------------------------------------------------------------
$cellDesc3 = $t.Cell(4, 4)
$desc3Text = $cellDesc3.Range.Text
$desc3Start = $cellDesc3.Range.Start
$oldStr2 = "Retrospectiva da Sprint-7 e início da Sprint-8"
$newStr2 = "Retrospectiva da Sprint-8 e início da Sprint-9"
$idx2 = $desc3Text.IndexOf($oldStr2)
if ($idx2 -ge 0) {
    $rngStart = $desc3Start + $idx2
    $rngEnd = $rngStart + $oldStr2.Length
    $d.Range($rngStart, $rngEnd).Text = $newStr2
}

# ---------------------------------------------------------------------------
# 4) Delete the whole "02" row (Nº 02, 27/10, 0,5, discussion about diagrams)
# ---------------------------------------------------------------------------
$t.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 5) Former row "03" (now row 3 after the deletion) becomes row "02":
#    Nº "03" -> "02", date "01/11" -> "08/11"
# ---------------------------------------------------------------------------
$cellNo3 = $t.Cell(3, 1)
$cellNo3.Range.Text = "02"

$cellDate3 = $t.Cell(3, 2)
$cellDate3.Range.Text = "08/11"

# ---------------------------------------------------------------------------
# 6) Column width fix: every description cell (column 4) width 6899 -> 6900
#    twips, i.e. 344.95pt -> 345pt
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $t.Cell($i, 4).Width = 345
}
